try {
    $x = [System.IO.Compression.ZipFile]::OpenRead("before.pptx")
    Write-Host "Opened" $x
} catch {
    Write-Host "ERR:" $_.Exception.Message
}
